# Apply the "experiment_data" cleanup edit described by the commit diff.
# All substantive changes are text edits inside Sheet1 (A1:H11): the
# Description column text is cleaned up (stray leading "- "/"'" characters
# removed, "." -> "," separators between the two model sentences), the
# "Diagonal Relocation-4" experiment name gets an underscore, and the
# "About the Agent Learning" header is wrapped in a <div> tag. The active
# selection also moves from I8 to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Experiment name cleanup (row 5, column A) ---------------------------
$ws.Range("A5").Value = "Diagonal_Relocation-4"

# --- Header cleanup (row 1, column D) -------------------------------------
$ws.Range("D1").Value = "<div style=`"width:290px`">About the Agent Learning</div>"

# --- Description column cleanup (column D, rows 2-5) ----------------------
$ws.Range("D2").Value = "Agent learns to stabalize it in the same position. Model - 0, Agent learns to slowly reach Goal . Model - 1"
$ws.Range("D3").Value = "Agent learns to stabalize it in the same position. Model - 0"
$ws.Range("D4").Value = "Agent learns to stabalize it in the same position. Model - 0 ,- Agent learns to slowly reach Goal . Model - 1"
$ws.Range("D5").Value = "Agent learns to reach goal at good time by rolling and attempts to stabalize it. Model - 0,  Agent learns to reach goal at good time by tapping and attempts to stabalize it.  . Model - 1"

# --- Move the active selection from I8 to I3 ------------------------------
$ws.Range("I3").Select()
